$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("C_16.1")

# Insert a new row above row 6 (shifts existing rows 6.. down by one)
$ws.Rows("6:6").Insert()

# Copy the formatting from the row that now matches the desired pattern for the
# new row (row 8, which carries the same style pair as the target row 6: s=10/s=11)
$ws.Range("B8:G8").Copy()
$ws.Range("B6:G6").PasteSpecial(-4122)

# Populate the new row 6 with the April 2025 data
$ws.Range("B6").Value = 2025
$ws.Range("C6").Value = "Abr."
$ws.Range("D6").Value = 9777
$ws.Range("E6").Value = 7179
$ws.Range("F6").Value = 14041
$ws.Range("G6").Value = 1452

# Correct the revised figures on the following three rows (now rows 7-9,
# formerly rows 6-8) that changed slightly with this update
$ws.Range("D7").Value = 10003
$ws.Range("E7").Value = 6771

$ws.Range("E8").Value = 5711

$ws.Range("D9").Value = 9110
$ws.Range("E9").Value = 5479

# Update the "Actualización" footnote (now on row 94) to reference April 2025
$ws.Range("B94").Value = "Actualización: Abril 2025."

# Grow the table / autofilter range to include the newly inserted row
$lo = $ws.ListObjects.Item("Tabla1")
$lo.Resize($ws.Range("B5:G93"))
